$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 488, shifting existing rows 488:504 down to 489:505.
$ws.Rows.Item(488).Insert()

# Populate the newly inserted row 488 with the new weekly data entry.
$ws.Range("A488").Value = 11
$ws.Range("B488").Value = 'Vega Monumental Concepción'
$ws.Range("C488").Value = 'Bíobío'
$ws.Range("D488").Value = 45041
$ws.Range("E488").Value = 8
$ws.Range("F488").Value = 100112002
$ws.Range("G488").Value = 'Pimiento'
$ws.Range("H488").Value = 'Morrón rojo'
$ws.Range("I488").Value = 'Primera'
$ws.Range("J488").Value = 190
$ws.Range("K488").Value = 10000
$ws.Range("L488").Value = 11000
$ws.Range("M488").Value = 10526
$ws.Range("N488").Value = '$/caja 18 kilos'
$ws.Range("O488").Value = 'Provincia de Limarí'
$ws.Range("P488").Value = 585
$ws.Range("Q488").Value = 18
$ws.Range("R488").Value = 'Hortaliza'
